$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 590408
$ws.Range("E2").Value = 18286
$ws.Range("F2").Value = 18286
$ws.Range("G2").Value = 12183
$ws.Range("H2").Value = 5014
$ws.Range("I2").Value = 3994
$ws.Range("J2").Value = 1020
$ws.Range("K2").Value = 370684
$ws.Range("L2").Value = 240773
$ws.Range("M2").Value = 129911
$ws.Range("N2").Value = 117194
$ws.Range("O2").Value = 12717
$ws.Range("P2").Value = 9042
$ws.Range("Q2").Value = 20292
$ws.Range("R2").Value = -23967
$ws.Range("S2").Value = 643
$ws.Range("T2").Value = 21659
$ws.Range("U2").Value = -1368
$ws.Range("V2").Value = 90024
$ws.Range("W2").Value = 3.1
$ws.Range("X2").Value = 0.85
$ws.Range("Y2").Value = 3.4
$ws.Range("Z2").Value = 1.38
$ws.Range("AA2").Value = 185.34
$ws.Range("AB2").Value = 1327.6
$ws.Range("AC2").Value = 2208
$ws.Range("AD2").Value = 26.76
$ws.Range("AE2").Value = 65084
$ws.Range("AF2").Value = 0.91
$ws.Range("AG2").Value = 400
$ws.Range("AH2").Value = 0.68
$ws.Range("AI2").Value = 18.25
$ws.Range("AJ2").Value = 163647814

# Row 3
$ws.Range("D3").Value = 565090
$ws.Range("E3").Value = 11923
$ws.Range("F3").Value = 11923
$ws.Range("G3").Value = 5931
$ws.Range("H3").Value = 2491
$ws.Range("I3").Value = 1244
$ws.Range("J3").Value = 1247
$ws.Range("K3").Value = 363139
$ws.Range("L3").Value = 233304
$ws.Range("M3").Value = 129835
$ws.Range("N3").Value = 116266
$ws.Range("O3").Value = 13569
$ws.Range("P3").Value = 9042
$ws.Range("Q3").Value = 26188
$ws.Range("R3").Value = -19332
$ws.Range("S3").Value = -1889
$ws.Range("T3").Value = 17473
$ws.Range("U3").Value = 8715
$ws.Range("V3").Value = 88272
$ws.Range("W3").Value = 2.11
$ws.Range("X3").Value = 0.44
$ws.Range("Y3").Value = 1.07
$ws.Range("Z3").Value = 0.68
$ws.Range("AA3").Value = 179.69
$ws.Range("AB3").Value = 1320.47
$ws.Range("AC3").Value = 688
$ws.Range("AD3").Value = 78.20999999999999
$ws.Range("AE3").Value = 64568
$ws.Range("AF3").Value = 0.83
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 0.74
$ws.Range("AI3").Value = 58.59
$ws.Range("AJ3").Value = 163647814

# Row 4
$ws.Range("D4").Value = 553670
$ws.Range("E4").Value = 13378
$ws.Range("F4").Value = 13378
$ws.Range("G4").Value = 7217
$ws.Range("H4").Value = 1263
$ws.Range("I4").Value = 769
$ws.Range("J4").Value = 494
$ws.Range("K4").Value = 378553
$ws.Range("L4").Value = 244985
$ws.Range("M4").Value = 133567
$ws.Range("N4").Value = 119871
$ws.Range("O4").Value = 13696
$ws.Range("P4").Value = 9042
$ws.Range("Q4").Value = 31580
$ws.Range("R4").Value = -23907
$ws.Range("S4").Value = -2788
$ws.Range("T4").Value = 20190
$ws.Range("U4").Value = 11390
$ws.Range("V4").Value = 86590
$ws.Range("W4").Value = 2.42
$ws.Range("X4").Value = 0.23
$ws.Range("Y4").Value = 0.65
$ws.Range("Z4").Value = 0.34
$ws.Range("AA4").Value = 183.42
$ws.Range("AB4").Value = 1344.52
$ws.Range("AC4").Value = 425
$ws.Range("AD4").Value = 121.37
$ws.Range("AE4").Value = 66571
$ws.Range("AF4").Value = 0.78
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 0.78
$ws.Range("AI4").Value = 94.81
$ws.Range("AJ4").Value = 163647814

# Row 5
$ws.Range("D5").Value = 613963
$ws.Range("E5").Value = 24685
$ws.Range("F5").Value = 24685
$ws.Range("G5").Value = 25581
$ws.Range("H5").Value = 18695
$ws.Range("I5").Value = 17258
$ws.Range("J5").Value = 1437
$ws.Range("K5").Value = 412210
$ws.Range("L5").Value = 265473
$ws.Range("M5").Value = 146737
$ws.Range("N5").Value = 132243
$ws.Range("O5").Value = 14494
$ws.Range("P5").Value = 9042
$ws.Range("Q5").Value = 21663
$ws.Range("R5").Value = -25829
$ws.Range("S5").Value = 8408
$ws.Range("T5").Value = 25755
$ws.Range("U5").Value = -4093
$ws.Range("V5").Value = 94505
$ws.Range("W5").Value = 4.02
$ws.Range("X5").Value = 3.04
$ws.Range("Y5").Value = 13.69
$ws.Range("Z5").Value = 4.73
$ws.Range("AA5").Value = 180.92
$ws.Range("AB5").Value = 1535.94
$ws.Range("AC5").Value = 9543
$ws.Range("AD5").Value = 11.11
$ws.Range("AE5").Value = 73441
$ws.Range("AF5").Value = 1.44
$ws.Range("AG5").Value = 400
$ws.Range("AH5").Value = 0.38
$ws.Range("AI5").Value = 4.22
$ws.Range("AJ5").Value = 163647814

# Row 6
$ws.Range("D6").Value = 613417
$ws.Range("E6").Value = 27033
$ws.Range("F6").Value = 27033
$ws.Range("G6").Value = 20086
$ws.Range("H6").Value = 14728
$ws.Range("I6").Value = 12401
$ws.Range("K6").Value = 443284
$ws.Range("L6").Value = 280215
$ws.Range("M6").Value = 163069
$ws.Range("N6").Value = 142533
$ws.Range("P6").Value = 9042
$ws.Range("Q6").Value = 45416
$ws.Range("R6").Value = -44203
$ws.Range("S6").Value = 8193
$ws.Range("T6").Value = 31665
$ws.Range("U6").Value = 13751
$ws.Range("V6").Value = 109012
$ws.Range("W6").Value = 4.41
$ws.Range("X6").Value = 2.4
$ws.Range("Y6").Value = 9.029999999999999
$ws.Range("Z6").Value = 3.44
$ws.Range("AA6").Value = 171.84
$ws.Range("AB6").Value = 1658.84
$ws.Range("AC6").Value = 6858
$ws.Range("AD6").Value = 9.08
$ws.Range("AE6").Value = 79156
$ws.Range("AF6").Value = 0.79
$ws.Range("AG6").Value = 750
$ws.Range("AH6").Value = 1.2
$ws.Range("AI6").Value = 10.96
$ws.Range("AJ6").Value = 163647814

# Row 7
$ws.Range("D7").Value = 626843
$ws.Range("E7").Value = 26124
$ws.Range("G7").Value = 14334
$ws.Range("H7").Value = 10635
$ws.Range("I7").Value = 8630
$ws.Range("K7").Value = 460878
$ws.Range("L7").Value = 288902
$ws.Range("M7").Value = 171975
$ws.Range("N7").Value = 150890
$ws.Range("P7").Value = 9041
$ws.Range("Q7").Value = 37859
$ws.Range("R7").Value = -31603
$ws.Range("S7").Value = -4758
$ws.Range("T7").Value = 25660
$ws.Range("U7").Value = 11436
$ws.Range("W7").Value = 4.17
$ws.Range("X7").Value = 1.7
$ws.Range("Y7").Value = 5.88
$ws.Range("Z7").Value = 2.35
$ws.Range("AA7").Value = 167.99
$ws.Range("AC7").Value = 4772
$ws.Range("AD7").Value = 14.9
$ws.Range("AE7").Value = 83797
$ws.Range("AF7").Value = 0.85
$ws.Range("AG7").Value = 796
$ws.Range("AH7").Value = 1.12
$ws.Range("AI7").Value = 15.1

# Row 8
$ws.Range("D8").Value = 653030
$ws.Range("E8").Value = 27204
$ws.Range("G8").Value = 18461
$ws.Range("H8").Value = 13469
$ws.Range("I8").Value = 12037
$ws.Range("K8").Value = 471403
$ws.Range("L8").Value = 294039
$ws.Range("M8").Value = 177364
$ws.Range("N8").Value = 154992
$ws.Range("P8").Value = 9041
$ws.Range("Q8").Value = 42354
$ws.Range("R8").Value = -29080
$ws.Range("S8").Value = -1883
$ws.Range("T8").Value = 24755
$ws.Range("U8").Value = 16655
$ws.Range("W8").Value = 4.17
$ws.Range("X8").Value = 2.06
$ws.Range("Y8").Value = 7.87
$ws.Range("Z8").Value = 2.89
$ws.Range("AA8").Value = 165.78
$ws.Range("AC8").Value = 6656
$ws.Range("AD8").Value = 9.9
$ws.Range("AE8").Value = 86075
$ws.Range("AF8").Value = 0.77
$ws.Range("AG8").Value = 811
$ws.Range("AH8").Value = 1.23
$ws.Range("AI8").Value = 11.03

# Row 9
$ws.Range("D9").Value = 678549
$ws.Range("E9").Value = 30852
$ws.Range("G9").Value = 22934
$ws.Range("H9").Value = 16768
$ws.Range("I9").Value = 15247
$ws.Range("K9").Value = 491921
$ws.Range("L9").Value = 299671
$ws.Range("M9").Value = 192249
$ws.Range("N9").Value = 169744
$ws.Range("P9").Value = 9041
$ws.Range("Q9").Value = 45847
$ws.Range("R9").Value = -29775
$ws.Range("S9").Value = -1550
$ws.Range("T9").Value = 24788
$ws.Range("U9").Value = 18888
$ws.Range("W9").Value = 4.55
$ws.Range("X9").Value = 2.47
$ws.Range("Y9").Value = 9.390000000000001
$ws.Range("Z9").Value = 3.48
$ws.Range("AA9").Value = 155.88
$ws.Range("AC9").Value = 8432
$ws.Range("AD9").Value = 7.82
$ws.Range("AE9").Value = 94267
$ws.Range("AF9").Value = 0.7
$ws.Range("AG9").Value = 824
$ws.Range("AH9").Value = 1.25
$ws.Range("AI9").Value = 8.84
